# __ADT_CONSTANTS__ update
#
# The "method id" constants get renamed from the old *_CMD_* naming
# convention to *_METHOD_*, and a brand new PRIORITY_QUEUE_METHOD_* block
# (5 rows) is inserted right after the STACK_METHOD_* block - pushing the
# LINKED_LIST_* and HASH_TABLE_* rows down from 37..60 to 42..65.
#
# Row 26 (STATUS_8033_FREE_SLOT_CORRUPTION) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rename the QUEUE_CMD_* / STACK_CMD_* constants in place -------
# (rows 27..36 keep their row position, only column A's text changes)
$renameInPlace = @{
    27 = "QUEUE_METHOD_00_NIL"
    28 = "QUEUE_METHOD_01_INIT"
    29 = "QUEUE_METHOD_10_ENQUEUE"
    30 = "QUEUE_METHOD_11_DEQUEUE"
    31 = "QUEUE_METHOD_12_PEEK"
    32 = "STACK_METHOD_00_NIL"
    33 = "STACK_METHOD_01_INIT"
    34 = "STACK_METHOD_10_PUSH"
    35 = "STACK_METHOD_11_POP"
    36 = "STACK_METHOD_12_PEEK"
}
foreach ($r in 27..36) {
    $ws.Cells.Item($r, 1).Value = $renameInPlace[$r]
}

# --- Step 2: make room for the new PRIORITY_QUEUE_METHOD_* block -----------
$ws.Rows.Item(37).Resize(5).Insert()

# --- Step 3: rename the (now shifted down by 5) LINKED_LIST_CMD_* / -------
#     HASH_TABLE_CMD_* constants in place (rows 42..65)
$renameShifted = @{
    42 = "LINKED_LIST_METHOD_00_NIL"
    43 = "LINKED_LIST_METHOD_01_INIT"
    44 = "LINKED_LIST_METHOD_10_ADD_FIRST"
    45 = "LINKED_LIST_METHOD_11_ADD_LAST"
    46 = "LINKED_LIST_METHOD_12_ADD_INDEX"
    47 = "LINKED_LIST_METHOD_13_ADD_FIRST_HIGH_PRIORITY"
    48 = "LINKED_LIST_METHOD_14_ADD_LAST_HIGH_PRIORITY"
    49 = "LINKED_LIST_METHOD_15_ADD_FIRST_LOW_PRIORITY"
    50 = "LINKED_LIST_METHOD_16_ADD_LAST_LOW_PRIORITY"
    51 = "LINKED_LIST_METHOD_20_REMOVE_FIRST"
    52 = "LINKED_LIST_METHOD_21_REMOVE_LAST"
    53 = "LINKED_LIST_METHOD_22_REMOVE_INDEX"
    54 = "LINKED_LIST_METHOD_23_REMOVE_FIRST_PRIORITY"
    55 = "LINKED_LIST_METHOD_24_REMOVE_LAST_PRIORITY"
    56 = "LINKED_LIST_METHOD_30_READ_FIRST"
    57 = "LINKED_LIST_METHOD_31_READ_LAST"
    58 = "LINKED_LIST_METHOD_32_READ_INDEX"
    59 = "LINKED_LIST_METHOD_33_READ_FIRST_PRIORITY"
    60 = "LINKED_LIST_METHOD_34_READ_LAST_PRIORITY"
    61 = "HASH_TABLE_METHOD_00_NIL"
    62 = "HASH_TABLE_METHOD_01_INIT"
    63 = "HASH_TABLE_METHOD_10_ADD"
    64 = "HASH_TABLE_METHOD_11_REMOVE"
    65 = "HASH_TABLE_METHOD_12_READ"
}
foreach ($r in 42..65) {
    $ws.Cells.Item($r, 1).Value = $renameShifted[$r]
}

# --- Step 4: fill in the 5 freshly inserted rows (37..41) with the new ----
#     PRIORITY_QUEUE_METHOD_* block
$path = "90_LIbrary\ADT\ADT"
$dataType = "int"

$priorityQueueRows = @(
    @("PRIORITY_QUEUE_METHOD_00_NIL",     0, " Method 00: Null "),
    @("PRIORITY_QUEUE_METHOD_01_INIT",    1, " Method 01: Initialize ADT "),
    @("PRIORITY_QUEUE_METHOD_10_ENQUEUE", 10, " Method 10: Add dato to last position "),
    @("PRIORITY_QUEUE_METHOD_11_DEQUEUE", 11, " Method 11: Remove data from first position "),
    @("PRIORITY_QUEUE_METHOD_12_PEEK",    12, " Method 12: Peek data from first position ")
)

$r = 37
foreach ($entry in $priorityQueueRows) {
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $path
    $ws.Cells.Item($r, 3).Value = $dataType
    $ws.Cells.Item($r, 4).Value = $entry[1]
    $ws.Cells.Item($r, 5).Value = $entry[2]
    $r++
}

$ws.Range("E6").Select()
